# "Generate Report for Handback"
# Fills in the handback columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) for the zh-cn and de-de localization sheets, and
# flips the Overview status from "Ready for handoff" to
# "Handed back: in sync with en-US". Also widens a few columns that now need
# to show the longer status / link text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$mdFileName = "68aeb95c-19ca-4db1-a4f6-7c5a94a4946d.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4597fd4b7efa8ac0d1d1221ef0a8f2eac4b60411/e2e/" + $mdFileName

# --- Overview sheet: status text + wider columns for it ---------------------
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("J2").Value = "68aeb95c-19ca-4db1-a4f6-7c5a94a4946d.380230b49e7c96d04f0098e0e5e9e7a06b1b4e47.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-28 12:59:41"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdFileName)

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("J2").Value = "68aeb95c-19ca-4db1-a4f6-7c5a94a4946d.380230b49e7c96d04f0098e0e5e9e7a06b1b4e47.de-de.xlf"
$dede.Range("K2").Value = "2016-08-28 12:59:48"

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdFileName)

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated."
